$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text content changes
$ws.Range("B2").Value = "rohan"
$ws.Range("L2").Value = "Reading ,Drawing"
$ws.Range("B3").Value = "mini"
$ws.Range("H3").Value = "abcd"
$ws.Range("L3").Value = "Reading ,Writing"

# Row heights: 18.75 -> 19.5 for header + both data rows
$ws.Rows("1:3").RowHeight = 19.5
